$d = $word.ActiveDocument

# Find the paragraph that contains "ReadMe schreiben" and insert a new
# list paragraph right after it, matching the same list style/numbering.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*ReadMe schreiben*") {
        $insertionRange = $p.Range
        $insertionRange.Collapse(0)  # wdCollapseEnd
        $insertionRange.InsertParagraphAfter()

        # Move the range to the newly created paragraph
        $newPara = $p.Next()
        $newPara.Range.Text = "Tabellen etc. als .tex Dokumente speichern, damit direkt bei Latex hochladen"

        # Match run-level font formatting
        $newPara.Range.Font.Name = "Cambria"

        break
    }
}
